$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.487532377243042
$ws.Range("B1").Value = 3.350951910018921
$ws.Range("C1").Value = 2.691777467727661
$ws.Range("D1").Value = 2.401623010635376
$ws.Range("E1").Value = 1.793492197990417
